$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Mercredi): add Fin (D5) and Temps total (E5)
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("D5").Value = 0.75
$ws.Range("E5").Value = "9h00"

# Row 6 (Jeudi): add Début (B6)
$ws.Range("B6").NumberFormat = "h:mm"
$ws.Range("B6").Value = 0.375

# Update active selection to B7
$ws.Range("B7").Select()
